# BattleShipEnemyDataExcel.xlsx edit
# 1. Plane Support 2. Ready to item table 3. add critical
#
# Adds two new columns ("CriticalRate" / "CriticalDamage") right before the
# existing "FireTime" column (old column P), shifting FireTime..the two
# trailing blank columns two places to the right. Fills in the new columns
# for the eight populated enemy rows, restores the custom column widths for
# the shifted columns plus the two new ones, bumps the header row height and
# updates the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1/2. Insert two columns at P:Q (shifts old P:U -> R:W) -----------------
$ws.Columns("P:Q").Insert()

# --- 3. New header labels ----------------------------------------------------
$ws.Range("P1").Value = "CriticalRate"
$ws.Range("Q1").Value = "CriticalDamage"

# New column data for the 8 populated enemy rows (rows 2-8)
$ws.Range("P2:P8").Value = 5
$ws.Range("Q2:Q8").Value = 5

# --- Column widths: restore custom widths on the two new columns ------------
$ws.Columns("P").ColumnWidth = 10
$ws.Columns("Q").ColumnWidth = 18.285714285714285

# --- Header row is now taller (wrapped 2-line headers) ----------------------
$ws.Rows(1).RowHeight = 26.25

# --- Saved selection moves to Q15 --------------------------------------------
$ws.Range("Q15").Select() | Out-Null

Write-Output "done"
